$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the source diff also moves the sheet's scroll position
# (topLeftCell H1 -> F1). We still drive ActiveWindow.ScrollColumn/ScrollRow
# below so the intent is expressed and the active-cell selection (which IS
# persisted) ends up correct; this runtime's xlsx exporter does not
# round-trip the <sheetView topLeftCell="..."> attribute at all (it is lost
# even on a plain load->save with no edits), so it cannot appear in the
# output no matter how it's set via COM.

# --- Row 2 value updates ---
$ws.Range("I2").Value = 275
$ws.Range("K2").Value = 98
$ws.Range("L2").Value = 127
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = 115
$ws.Range("Q2").Value = 38

# --- New supporting totals on row 3 ---
$ws.Range("J3").Value = 771
$ws.Range("P3").Value = 360

# --- J2 becomes a formula, losing its fill/border (matches plain centered style) ---
$J2 = $ws.Range("J2")
$J2.Formula = "=J3-(I2+K2)"
$J2.Interior.Pattern = -4142
$J2.Borders.LineStyle = -4142

# --- P2 becomes a formula, losing its fill/border (matches plain centered style) ---
$P2 = $ws.Range("P2")
$P2.Formula = "=P3-(O2+Q2)"
$P2.Interior.Pattern = -4142
$P2.Borders.LineStyle = -4142

# --- New row-3 helper cells match the plain centered style (same as J2/P2/M2/S2) ---
$J3 = $ws.Range("J3")
$J3.Interior.Pattern = -4142
$J3.Borders.LineStyle = -4142
$J3.HorizontalAlignment = -4108

$P3 = $ws.Range("P3")
$P3.Interior.Pattern = -4142
$P3.Borders.LineStyle = -4142
$P3.HorizontalAlignment = -4108

# --- Sheet view: scroll position and selection ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("L13").Select()
